$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the updated crypto price/volume values (and the two swapped rows, FraxShare <-> VeChain).
# Price-column cells whose new text is numeric-looking get an explicit text NumberFormat first,
# so Excel stores them verbatim (preserves trailing zeros / avoids scientific notation) instead of
# silently converting the literal into a floating-point Number.
$ws.Range("D2").Value = "29.661.70"
$ws.Range("E2").Value = "  -1.79%  "
$ws.Range("D3").Value = "2.005.68"
$ws.Range("E3").Value = "  -3.93%  "
$ws.Range("E4").Value = "  +0.57%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "331.01"
$ws.Range("E5").Value = "  -3.58%  "
$ws.Range("E6").Value = "  +0.49%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5020"
$ws.Range("E7").Value = "  -3.58%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4257"
$ws.Range("E8").Value = "  -3.11%  "
$ws.Range("E9").Value = "  +0.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08974"
$ws.Range("E10").Value = "  -3.31%  "
$ws.Range("E11").Value = "  -3.82%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.42"
$ws.Range("E12").Value = "  -4.96%  "
$ws.Range("D13").Value = "2.047.27"
$ws.Range("E13").Value = "  -3.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.112"
$ws.Range("E14").Value = "  -6.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.508"
$ws.Range("E15").Value = "  -5.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.013"
$ws.Range("E16").Value = "  +0.46%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "94.49"
$ws.Range("E17").Value = "  -6.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001116"
$ws.Range("E18").Value = "  -3.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06671"
$ws.Range("E19").Value = "  -0.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.82"
$ws.Range("E20").Value = "  -6.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.982"
$ws.Range("E22").Value = "  -5.35%  "
$ws.Range("D23").Value = "29.666.11"
$ws.Range("E23").Value = "  -1.91%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.05"
$ws.Range("E24").Value = "  -3.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.282"
$ws.Range("E25").Value = "  -0.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.46"
$ws.Range("E26").Value = "  -1.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.79"
$ws.Range("E27").Value = "  -4.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.370"
$ws.Range("E28").Value = "  -4.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.314"
$ws.Range("E29").Value = "  -7.47%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "129.15"
$ws.Range("E30").Value = "  -2.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.062"
$ws.Range("E31").Value = "  -5.63%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09969"
$ws.Range("E32").Value = "  -4.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.573"
$ws.Range("E34").Value = "  -5.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.816"
$ws.Range("E35").Value = "  -1.50%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02478"
$ws.Range("E36").Value = "  -5.26%  "
$ws.Range("B37").Value = "FraxShare"
$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.440"
$ws.Range("E37").Value = "  -7.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.318"
$ws.Range("E38").Value = "  -2.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06374"
$ws.Range("E39").Value = "  -5.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6595"
$ws.Range("E40").Value = "  -5.21%  "
$ws.Range("E41").Value = "  -6.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2059"
$ws.Range("E42").Value = "  -6.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.011"
$ws.Range("E43").Value = "  +0.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6361"
$ws.Range("E44").Value = "  -6.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.60"
$ws.Range("E45").Value = "  -5.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.213"
$ws.Range("E46").Value = "  -5.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.320"
$ws.Range("E47").Value = "  -3.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.529"
$ws.Range("E48").Value = "  -3.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000340"
$ws.Range("E49").Value = "  -2.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07003"
$ws.Range("E50").Value = "  -2.94%  "
$ws.Range("E51").Value = "  -6.48%  "
